$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update status-check timestamp in F1
$ws.Range("F1").Value = "Last status check on: 24.02.2022 10:00"

# Update row 6 (Shell Olomoucka): Delta Cena becomes a numeric value,
# and Old Datum becomes a real date/time value formatted like the other rows.
$ws.Range("D6").Value = 0.39
$ws.Range("E6").Value = 44616.40800925926
$ws.Range("E6").NumberFormat = "YYYY-MM-DD HH:MM:SS"
